# [F2] Change expected filename in test-data
# - Update "Expected filename on server" for F2-1-2 and F2-1-3 from *.txt to *.php
# - Add a new "Result" column (F) with Passed/Failed status for each test row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update E3 (F2-1-2) and E4 (F2-1-3): change file extension from .txt to .php ---
$e3 = $ws.Range("E3").Value()
$e3base = $e3.Substring(0, $e3.Length - 4)
$ws.Range("E3").Value = $e3base + ".php"

$e4 = $ws.Range("E4").Value()
$e4base = $e4.Substring(0, $e4.Length - 4)
$ws.Range("E4").Value = $e4base + ".php"

# --- Add new column F with pass/fail results for each data row ---
$ws.Range("F2").Value = "Passed"
$ws.Range("F3").Value = "Failed"
$ws.Range("F4").Value = "Failed"
$ws.Range("F5").Value = "Passed"
$ws.Range("F6").Value = "Passed"
$ws.Range("F7").Value = "Passed"
$ws.Range("F8").Value = "Passed"

# --- Restore selection state ---
$ws.Range("K10").Select()
